# Lower-case the E1:G1 identifier header cells on the DataSubstrate,
# DataTopic and Organization sheets (and their "1"-suffixed duplicates).

$wb = $excel.ActiveWorkbook

$sheetsEdamGroup = @("DataSubstrate", "DataTopic", "DataSubstrate1", "DataTopic1")
foreach ($name in $sheetsEdamGroup) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("E1").Value = "edam_id"
    $ws.Range("F1").Value = "mesh_id"
    $ws.Range("G1").Value = "ncit_id"
}

$sheetsOrgGroup = @("Organization", "Organization1")
foreach ($name in $sheetsOrgGroup) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("E1").Value = "ror_id"
    $ws.Range("F1").Value = "wikidata_id"
    $ws.Range("G1").Value = "url"
}
